$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("G2").Value = 30.93275633333333
$ws.Range("H2").Value = 92.798269
$ws.Range("I2").Value = 0.6015903584115317
$ws.Range("J2").Value = 0.6015903584115317
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.125191
$ws.Range("N2").Value = 0.375573
$ws.Range("Q2").Value = 3.872502698126333
$ws.Range("R2").Value = 34.852524283137
$ws.Range("S2").Value = 0.6015903584115317
$ws.Range("T2").Value = 0.6015903584115317

# Row 3
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("I3").Value = 0.002688776579266707
$ws.Range("J3").Value = 0.002688776579266707
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.125191
$ws.Range("N3").Value = 0.375573
$ws.Range("Q3").Value = 0.01730794786233333
$ws.Range("R3").Value = 0.155771530761
$ws.Range("S3").Value = 0.002688776579266707
$ws.Range("T3").Value = 0.002688776579266707

# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 13.65672433333333
$ws.Range("H4").Value = 40.970173
$ws.Range("I4").Value = 0.2656004398018724
$ws.Range("J4").Value = 0.2656004398018724
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.125191
$ws.Range("N4").Value = 0.375573
$ws.Range("Q4").Value = 1.709698976014334
$ws.Range("R4").Value = 15.387290784129
$ws.Range("S4").Value = 0.2656004398018724
$ws.Range("T4").Value = 0.2656004398018724

# Row 5
$ws.Range("D5").Value = "Inflammatory-Mac"
$ws.Range("G5").Value = 0.1180373333333333
$ws.Range("H5").Value = 0.354112
$ws.Range("I5").Value = 0.002295628650118725
$ws.Range("J5").Value = 0.002295628650118725
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.125191
$ws.Range("N5").Value = 0.375573
$ws.Range("Q5").Value = 0.01477721179733333
$ws.Range("R5").Value = 0.132994906176
$ws.Range("S5").Value = 0.002295628650118725
$ws.Range("T5").Value = 0.002295628650118725

# Row 6
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("G6").Value = 6.572534333333333
$ws.Range("H6").Value = 19.717603
$ws.Range("I6").Value = 0.1278247965572105
$ws.Range("J6").Value = 0.1278247965572105
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.125191
$ws.Range("N6").Value = 0.375573
$ws.Range("Q6").Value = 0.8228221457243332
$ws.Range("R6").Value = 7.405399311518998
$ws.Range("S6").Value = 0.1278247965572105
$ws.Range("T6").Value = 0.1278247965572105
